$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Language codes already used elsewhere in the sheet, in the order the
# new template-type blocks were appended with.
$langs = @("eng", "ara", "fra", "hin", "kan", "tam")

# Each new template type: (code, description)
$templates = @(
    @("tnc-order-a-physical-card", "Order a physical card"),
    @("tnc-share-cred-with-partner", "Share your credential with a partner"),
    @("tnc-update-demo", "Update demographic data")
)

$row = 914
foreach ($tpl in $templates) {
    $code = $tpl[0]
    $descr = $tpl[1]
    foreach ($lang in $langs) {
        $ws.Range("A$row").Value = $lang
        $ws.Range("B$row").Value = $code
        $ws.Range("C$row").Value = $descr
        # Re-use an existing "TRUE" text cell so the value lands as the
        # shared text string "TRUE" (matching the rest of column D)
        # instead of a native boolean.
        $ws.Range("D2").Copy($ws.Range("D$row"))
        $row = $row + 1
    }
}
